$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Read the whole country table (A4:H216) into memory so we can reorder it
#    (Catar moves up next to Indonesia, Maldivas moves up next to Bermudas)
#    and patch the handful of cells whose figures were refreshed.
# ---------------------------------------------------------------------------
$nRows = 213
$nCols = 8
$rng = $ws.Range("A4:H216")
$data = $rng.Value2

$rows = New-Object System.Collections.ArrayList
for ($i = 1; $i -le $nRows; $i++) {
    $row = New-Object 'object[]' $nCols
    for ($j = 1; $j -le $nCols; $j++) {
        $row[$j - 1] = $data[$i, $j]
    }
    $rows.Add($row) | Out-Null
}

# ---------------------------------------------------------------------------
# 2) Move "Catar" so it sits right after "Indonesia" (ahead of "Noruega").
# ---------------------------------------------------------------------------
$moveIdx = -1
for ($k = 0; $k -lt $rows.Count; $k++) {
    if ($rows[$k][0] -eq "Catar") { $moveIdx = $k }
}
$rowData = $rows[$moveIdx]
$rows.RemoveAt($moveIdx)
$afterIdx = -1
for ($k = 0; $k -lt $rows.Count; $k++) {
    if ($rows[$k][0] -eq "Indonesia") { $afterIdx = $k }
}
$insertAt = $afterIdx + 1
$rows.Insert($insertAt, $rowData)

# ---------------------------------------------------------------------------
# 3) Move "Maldivas" so it sits right after "Bermudas" (ahead of "Monaco").
# ---------------------------------------------------------------------------
$moveIdx2 = -1
for ($k = 0; $k -lt $rows.Count; $k++) {
    if ($rows[$k][0] -eq "Maldivas") { $moveIdx2 = $k }
}
$rowData2 = $rows[$moveIdx2]
$rows.RemoveAt($moveIdx2)
$afterIdx2 = -1
for ($k = 0; $k -lt $rows.Count; $k++) {
    if ($rows[$k][0] -eq "Bermudas") { $afterIdx2 = $k }
}
$insertAt2 = $afterIdx2 + 1
$rows.Insert($insertAt2, $rowData2)

# ---------------------------------------------------------------------------
# 4) Refresh the figures for the countries whose counts changed.
#    Columns: B=Casos totales C=Nuevos casos D=Casos activos E=Recuperados
#             F=Casos criticos G=Muertes hoy H=Muertes
# ---------------------------------------------------------------------------
for ($k = 0; $k -lt $rows.Count; $k++) {
    $name = $rows[$k][0]
    if ($name -eq "Iran") {
        $rows[$k][1] = 87026
        $rows[$k][2] = 1030
        $rows[$k][3] = 64843
        $rows[$k][4] = 16702
        $rows[$k][5] = 3105
        $rows[$k][6] = 90
        $rows[$k][7] = 5481
    }
    elseif ($name -eq "Emiratos Arabes Unidos") {
        $rows[$k][1] = 8756
        $rows[$k][2] = 518
        $rows[$k][3] = 1637
        $rows[$k][4] = 7063
        $rows[$k][6] = 4
        $rows[$k][7] = 56
    }
    elseif ($name -eq "Catar") {
        $rows[$k][1] = 7764
        $rows[$k][2] = 623
        $rows[$k][3] = 750
        $rows[$k][4] = 7004
        $rows[$k][5] = 37
        $rows[$k][6] = 0
        $rows[$k][7] = 10
    }
    elseif ($name -eq "Malta") {
        $rows[$k][1] = 445
        $rows[$k][2] = 1
        $rows[$k][3] = 204
        $rows[$k][4] = 238
    }
    elseif ($name -eq "Madagascar") {
        $rows[$k][3] = 58
        $rows[$k][4] = 63
    }
    elseif ($name -eq "Maldivas") {
        $rows[$k][1] = 94
        $rows[$k][2] = 8
        $rows[$k][3] = 16
        $rows[$k][4] = 78
        $rows[$k][7] = 0
    }
    elseif ($name -eq "Nepal") {
        $rows[$k][1] = 47
        $rows[$k][2] = 2
        $rows[$k][3] = 9
    }
}

# ---------------------------------------------------------------------------
# 5) Write the reordered/updated table back out in one shot.
# ---------------------------------------------------------------------------
$out = New-Object 'object[,]' $nRows,$nCols
for ($i = 1; $i -le $nRows; $i++) {
    for ($j = 1; $j -le $nCols; $j++) {
        $out[$i - 1, $j - 1] = $rows[$i - 1][$j - 1]
    }
}
$ws.Range("A4:H216").Value2 = $out

# ---------------------------------------------------------------------------
# 6) Bump the "last updated" timestamp banner in A1.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value2 = "Datos actualizados a 23 de Abril de 2020 a las 12:52"
